$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = "System, backup@backdoor.com, system"
$ws.Range("G3").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G4").Value2 = "System, backup@backdoor.com"
$ws.Range("G5").Value2 = "System, backup@backdoor.com"
$ws.Range("G6").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G10").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G11").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G12").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G13").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G14").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G15").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G17").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G18").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G19").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G20").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G29").Value2 = "System, backup@backdoor.com, system"
$ws.Range("G30").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G31").Value2 = "System, backup@backdoor.com"
$ws.Range("G32").Value2 = "System, backup@backdoor.com"
$ws.Range("G33").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G37").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G38").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G39").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G40").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G41").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G42").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G44").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G45").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G46").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G47").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G56").Value2 = "System, backup@backdoor.com, system"
$ws.Range("G57").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G58").Value2 = "System, backup@backdoor.com"
$ws.Range("G59").Value2 = "System, backup@backdoor.com"
$ws.Range("G60").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G64").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G65").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G66").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G67").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G68").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G69").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G71").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G72").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G73").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G74").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G83").Value2 = "System, backup@backdoor.com"
$ws.Range("G84").Value2 = "System, backup@backdoor.com"
$ws.Range("G85").Value2 = "System, backup@backdoor.com"
$ws.Range("G86").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G87").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G88").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G89").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G90").Value2 = "admin@admin.com, dnasr281@gmail.com"
$ws.Range("G93").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G95").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G96").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G97").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G99").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G109").Value2 = "System, backup@backdoor.com"
$ws.Range("G110").Value2 = "System, backup@backdoor.com"
$ws.Range("G111").Value2 = "System, backup@backdoor.com"
$ws.Range("G112").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G113").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G114").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G115").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G116").Value2 = "admin@admin.com, dnasr281@gmail.com"
$ws.Range("G119").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G121").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G122").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G123").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G125").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G135").Value2 = "System, backup@backdoor.com"
$ws.Range("G136").Value2 = "System, backup@backdoor.com"
$ws.Range("G137").Value2 = "System, backup@backdoor.com"
$ws.Range("G138").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G139").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G140").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G141").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G142").Value2 = "admin@admin.com, dnasr281@gmail.com"
$ws.Range("G145").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G147").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G148").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G149").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G151").Value2 = "System, dnasr281@gmail.com"
